# SZE_arajanlat_sablon.docx - template placeholder touch-up
#
# Actual semantic content changes applied by this commit:
#   1. The "Nettó ár" column placeholder {unit_price_fmt} -> {sum_total_fmt}
#   2. The "Áfa" column literal 0%                         -> {afa}
#   3. The "Bruttó ár" column placeholder {line_total_fmt}  -> {sum_gross_fmt}
#      (the trailing {/items} loop-close tag is left untouched)
#   4. The hidden _GoBack bookmark that was straddling the
#      "{helyszin}" placeholder is dropped (re-typing the placeholder
#      text removes the stale bookmark while leaving the visible text,
#      "{helyszin}", and its bold formatting unchanged).
#
# Everything else in the source diff (extra proofErr spell-check tags
# around already-correct words, additional namespace declarations,
# wp14:anchorId/editId attributes, bookmark renumbering, style id
# aliasing, etc.) is cosmetic fallout from re-saving the package with a
# newer Word build and carries no content/formatting change, so it is
# not re-created here.

$d = $word.ActiveDocument

# 1) Nettó ár (net price) placeholder rename.
$d.Content.Find.Execute(
    "{unit_price_fmt}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{sum_total_fmt}", 2) | Out-Null

# 2) Áfa (VAT) column: literal "0%" becomes a real placeholder.
$d.Content.Find.Execute(
    "0%", $true, $false, $false, $false, $false,
    $true, 1, $false, "{afa}", 2) | Out-Null

# 3) Bruttó ár (gross price) placeholder rename; keep "{/items}" intact.
$d.Content.Find.Execute(
    "{line_total_fmt}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{sum_gross_fmt}", 2) | Out-Null

# 4) Re-type "{helyszin}" verbatim so Word drops the leftover hidden
#    _GoBack bookmark that used to sit between "{helyszin" and "}".
$d.Content.Find.Execute(
    "{helyszin}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{helyszin}", 2) | Out-Null
